$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.101.78"
$ws.Range("E2").Value = "  -2.52%  "
$ws.Range("D3").Value = "2.232.17"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.41"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -5.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("D14").Value = "2.568.52"
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.848"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.50%  "
$ws.Range("D17").Value = "2.230.27"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").Value = "41.924.90"
$ws.Range("E18").Value = "  -2.83%  "
$ws.Range("D19").Value = "0.0₃0971"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  -5.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.49%  "
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0821"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.118"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0298"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "111.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.200"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0989"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  -5.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -15.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.46%  "
